$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dCell = $ws.Range("D2")
$dCell.NumberFormat = "@"
$dCell.Value = '28.571.66'
$dCell.Style = "Normal"
$ws.Range("E2").Value = '  +0.69%  '
$dCell = $ws.Range("D3")
$dCell.NumberFormat = "@"
$dCell.Value = '1.873.73'
$dCell.Style = "Normal"
$ws.Range("E3").Value = '  +0.37%  '
$dCell = $ws.Range("D4")
$dCell.NumberFormat = "@"
$dCell.Value = '1.006'
$dCell.Style = "Normal"
$ws.Range("E4").Value = '  -1.22%  '
$dCell = $ws.Range("D5")
$dCell.NumberFormat = "@"
$dCell.Value = '314.53'
$dCell.Style = "Normal"
$ws.Range("E5").Value = '  -0.76%  '
$dCell = $ws.Range("D6")
$dCell.NumberFormat = "@"
$dCell.Value = '1.005'
$dCell.Style = "Normal"
$ws.Range("E6").Value = '  -1.22%  '
$dCell = $ws.Range("D7")
$dCell.NumberFormat = "@"
$dCell.Value = '0.5079'
$dCell.Style = "Normal"
$ws.Range("E7").Value = '  -0.56%  '
$dCell = $ws.Range("D8")
$dCell.NumberFormat = "@"
$dCell.Value = '0.3910'
$dCell.Style = "Normal"
$ws.Range("E8").Value = '  -1.15%  '
$dCell = $ws.Range("D9")
$dCell.NumberFormat = "@"
$dCell.Value = '0.08346'
$dCell.Style = "Normal"
$ws.Range("E9").Value = '  -1.76%  '
$dCell = $ws.Range("D10")
$dCell.NumberFormat = "@"
$dCell.Value = '42.34'
$dCell.Style = "Normal"
$ws.Range("E10").Value = '  +1.06%  '
$dCell = $ws.Range("D11")
$dCell.NumberFormat = "@"
$dCell.Value = '1.106'
$dCell.Style = "Normal"
$ws.Range("E11").Value = '  -0.26%  '
$dCell = $ws.Range("D12")
$dCell.NumberFormat = "@"
$dCell.Value = '6.194'
$dCell.Style = "Normal"
$ws.Range("E12").Value = '  -0.72%  '
$dCell = $ws.Range("D13")
$dCell.NumberFormat = "@"
$dCell.Value = '1.869.14'
$dCell.Style = "Normal"
$ws.Range("E13").Value = '  +3.86%  '
$dCell = $ws.Range("D14")
$dCell.NumberFormat = "@"
$dCell.Value = '20.35'
$dCell.Style = "Normal"
$ws.Range("E14").Value = '  -0.30%  '
$dCell = $ws.Range("D15")
$dCell.NumberFormat = "@"
$dCell.Value = '7.260'
$dCell.Style = "Normal"
$ws.Range("E15").Value = '  +0.91%  '
$dCell = $ws.Range("D16")
$dCell.NumberFormat = "@"
$dCell.Value = '1.006'
$dCell.Style = "Normal"
$ws.Range("E16").Value = '  -1.17%  '
$dCell = $ws.Range("D17")
$dCell.NumberFormat = "@"
$dCell.Value = '93.31'
$dCell.Style = "Normal"
$ws.Range("E17").Value = '  +3.11%  '
$dCell = $ws.Range("D18")
$dCell.NumberFormat = "@"
$dCell.Value = '0.00001099'
$dCell.Style = "Normal"
$ws.Range("E18").Value = '  -1.46%  '
$dCell = $ws.Range("D19")
$dCell.NumberFormat = "@"
$dCell.Value = '0.06715'
$dCell.Style = "Normal"
$ws.Range("E19").Value = '  -0.56%  '
$dCell = $ws.Range("D20")
$dCell.NumberFormat = "@"
$dCell.Value = '17.64'
$dCell.Style = "Normal"
$ws.Range("E20").Value = '  -0.15%  '
$dCell = $ws.Range("D21")
$dCell.NumberFormat = "@"
$dCell.Value = '1.004'
$dCell.Style = "Normal"
$ws.Range("E21").Value = '  -1.26%  '
$dCell = $ws.Range("D22")
$dCell.NumberFormat = "@"
$dCell.Value = '5.935'
$dCell.Style = "Normal"
$ws.Range("E22").Value = '  -0.20%  '
$dCell = $ws.Range("D23")
$dCell.NumberFormat = "@"
$dCell.Value = '28.596.06'
$dCell.Style = "Normal"
$ws.Range("E23").Value = '  +0.78%  '
$dCell = $ws.Range("D24")
$dCell.NumberFormat = "@"
$dCell.Value = '11.09'
$dCell.Style = "Normal"
$ws.Range("E24").Value = '  -0.41%  '
$dCell = $ws.Range("D25")
$dCell.NumberFormat = "@"
$dCell.Value = '2.192'
$dCell.Style = "Normal"
$ws.Range("E25").Value = '  -3.93%  '
$dCell = $ws.Range("D26")
$dCell.NumberFormat = "@"
$dCell.Value = '2.082.95'
$dCell.Style = "Normal"
$ws.Range("E26").Value = '  +3.57%  '
$dCell = $ws.Range("D27")
$dCell.NumberFormat = "@"
$dCell.Value = '157.98'
$dCell.Style = "Normal"
$ws.Range("E27").Value = '  -2.47%  '
$dCell = $ws.Range("D28")
$dCell.NumberFormat = "@"
$dCell.Value = '20.56'
$dCell.Style = "Normal"
$ws.Range("E28").Value = '  -0.89%  '
$dCell = $ws.Range("D29")
$dCell.NumberFormat = "@"
$dCell.Value = '2.423'
$dCell.Style = "Normal"
$ws.Range("E29").Value = '  +3.03%  '
$dCell = $ws.Range("D30")
$dCell.NumberFormat = "@"
$dCell.Value = '125.92'
$dCell.Style = "Normal"
$ws.Range("E30").Value = '  -0.96%  '
$dCell = $ws.Range("D31")
$dCell.NumberFormat = "@"
$dCell.Value = '0.1037'
$dCell.Style = "Normal"
$ws.Range("E31").Value = '  -1.20%  '
$dCell = $ws.Range("D32")
$dCell.NumberFormat = "@"
$dCell.Value = '1.045'
$dCell.Style = "Normal"
$ws.Range("E32").Value = '  +1.16%  '
$dCell = $ws.Range("D33")
$dCell.NumberFormat = "@"
$dCell.Value = '5.775'
$dCell.Style = "Normal"
$ws.Range("E33").Value = '  -0.07%  '
$dCell = $ws.Range("D34")
$dCell.NumberFormat = "@"
$dCell.Value = '3.639'
$dCell.Style = "Normal"
$ws.Range("E34").Value = '  -0.28%  '
$dCell = $ws.Range("D35")
$dCell.NumberFormat = "@"
$dCell.Value = '0.02452'
$dCell.Style = "Normal"
$ws.Range("E35").Value = '  +1.25%  '
$dCell = $ws.Range("D36")
$dCell.NumberFormat = "@"
$dCell.Value = '0.06554'
$dCell.Style = "Normal"
$ws.Range("E36").Value = '  +1.45%  '
$dCell = $ws.Range("D37")
$dCell.NumberFormat = "@"
$dCell.Value = '9.014'
$dCell.Style = "Normal"
$ws.Range("E37").Value = '  +1.98%  '
$dCell = $ws.Range("D38")
$dCell.NumberFormat = "@"
$dCell.Value = '0.2164'
$dCell.Style = "Normal"
$ws.Range("E38").Value = '  -0.87%  '
$dCell = $ws.Range("D39")
$dCell.NumberFormat = "@"
$dCell.Value = '5.030'
$dCell.Style = "Normal"
$ws.Range("E39").Value = '  +0.78%  '
$dCell = $ws.Range("D40")
$dCell.NumberFormat = "@"
$dCell.Value = '1.187'
$dCell.Style = "Normal"
$ws.Range("E40").Value = '  +0.78%  '
$dCell = $ws.Range("D41")
$dCell.NumberFormat = "@"
$dCell.Value = '1.238'
$dCell.Style = "Normal"
$ws.Range("E41").Value = '  -1.91%  '
$dCell = $ws.Range("D42")
$dCell.NumberFormat = "@"
$dCell.Value = '0.6376'
$dCell.Style = "Normal"
$ws.Range("E42").Value = '  +0.01%  '
$dCell = $ws.Range("D43")
$dCell.NumberFormat = "@"
$dCell.Value = '11.14'
$dCell.Style = "Normal"
$ws.Range("E43").Value = '  -1.00%  '
$dCell = $ws.Range("D44")
$dCell.NumberFormat = "@"
$dCell.Value = '1.005'
$dCell.Style = "Normal"
$ws.Range("E44").Value = '  -0.86%  '
$dCell = $ws.Range("D45")
$dCell.NumberFormat = "@"
$dCell.Value = '0.5985'
$dCell.Style = "Normal"
$ws.Range("E45").Value = '  -0.52%  '
$dCell = $ws.Range("D46")
$dCell.NumberFormat = "@"
$dCell.Value = '13.00'
$dCell.Style = "Normal"
$ws.Range("E46").Value = '  -0.01%  '
$dCell = $ws.Range("D47")
$dCell.NumberFormat = "@"
$dCell.Value = '3.675'
$dCell.Style = "Normal"
$ws.Range("E47").Value = '  -0.84%  '
$dCell = $ws.Range("D48")
$dCell.NumberFormat = "@"
$dCell.Value = '2.004'
$dCell.Style = "Normal"
$ws.Range("E48").Value = '  +0.84%  '
$dCell = $ws.Range("D49")
$dCell.NumberFormat = "@"
$dCell.Value = '1.220'
$dCell.Style = "Normal"
$ws.Range("E49").Value = '  +1.45%  '
$dCell = $ws.Range("D50")
$dCell.NumberFormat = "@"
$dCell.Value = '122.36'
$dCell.Style = "Normal"
$ws.Range("E50").Value = '  +1.22%  '
$dCell = $ws.Range("D51")
$dCell.NumberFormat = "@"
$dCell.Value = '1.186'
$dCell.Style = "Normal"
$ws.Range("E51").Value = '  -2.26%  '
